$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh: update Price (D) and Volume(1h) (E) columns.
# Some Price values are plain decimal numbers (e.g. "586.84"); the source sheet
# stores these as plain text, so we briefly force a text number format before
# assignment (otherwise Excel auto-converts them into numeric cells), then clear
# the formatting again so the cell keeps its original (unformatted) style.

$ws.Range("D2").Value = "63.126.21"
$ws.Range("E2").Value = "  -0.03%  "

$ws.Range("D3").Value = "3.052.81"
$ws.Range("E3").Value = "  -0.13%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.84"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.22"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.52%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  -1.78%  "

$ws.Range("D9").Value = "3.053.66"
$ws.Range("E9").Value = "  -0.23%  "

$ws.Range("E10").Value = "  -1.44%  "

$ws.Range("E11").Value = "  +0.11%  "

$ws.Range("E12").Value = "  -2.58%  "

$ws.Range("E13").Value = "  -2.18%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.23"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.92%  "

$ws.Range("E15").Value = "  +1.87%  "

$ws.Range("D16").Value = "3.555.37"
$ws.Range("E16").Value = "  -0.26%  "

$ws.Range("E17").Value = "  -0.62%  "

$ws.Range("D18").Value = "63.082.05"
$ws.Range("E18").Value = "  -0.17%  "

$ws.Range("D19").Value = "3.050.91"
$ws.Range("E19").Value = "  -0.31%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "476.51"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.70%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.27"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.46%  "

$ws.Range("E22").Value = "  -1.49%  "

$ws.Range("E23").Value = "  -0.13%  "

$ws.Range("E24").Value = "  +1.69%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.30"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.51%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.68"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.61"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +6.45%  "

$ws.Range("E28").Value = "  +0.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.33"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.91%  "

$ws.Range("E30").Value = "  +0.08%  "

$ws.Range("E31").Value = "  -0.10%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.19"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.21%  "

$ws.Range("E33").Value = "  +1.75%  "

$ws.Range("E34").Value = "  -2.17%  "

$ws.Range("E35").Value = "  +1.33%  "

$ws.Range("D36").Value = "0.0₃0817"
$ws.Range("E36").Value = "  -2.63%  "

$ws.Range("E37").Value = "  -2.80%  "

$ws.Range("E38").Value = "  +0.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.89"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.31%  "

$ws.Range("E40").Value = "  -0.52%  "

$ws.Range("E41").Value = "  +0.28%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "433.20"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.287"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.06%  "

$ws.Range("E44").Value = "  +3.12%  "

$ws.Range("E45").Value = "  -0.12%  "

$ws.Range("D46").Value = "2.827.92"
$ws.Range("E46").Value = "  +1.25%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.39"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.70%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.18"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.20"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.70%  "

$ws.Range("E51").Value = "  -1.49%  "
